# End of chapter 3 - append a new "embedding size" term to the list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next empty row in column A and append the new term there.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1
$ws.Cells.Item($newRow, 1).Value = "embedding size"

# Restore the on-screen scroll position / selection as left by the author.
$excel.ActiveWindow.ScrollRow = 52
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("A11").Select()
